$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateColStyleSource = $ws.Range("A374")

# Row 375
$dateColStyleSource.Copy()
$ws.Range("A375").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_375 = @(44449,1,0,0,6,5,0,4,0,3,0,1,0,3,0,1,0,4,1,3,6,44,1,1,3,1,1,0,1,1,2,1,6,0,0,0,1,3,2,0,6,114,0,0,0,0,1,0,0,1)
for ($i = 0; $i -lt $rowVals_375.Length; $i++) {
    $ws.Cells.Item(375, $i + 1).Value = $rowVals_375[$i]
}

# Row 376
$dateColStyleSource.Copy()
$ws.Range("A376").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_376 = @(44450,1,2,0,24,6,0,1,0,1,0,0,1,0,0,0,0,1,0,0,0,22,1,2,18,5,0,1,4,0,1,0,0,1,0,0,4,0,3,0,5,106,1,0,0,0,1,0,0,0)
for ($i = 0; $i -lt $rowVals_376.Length; $i++) {
    $ws.Cells.Item(376, $i + 1).Value = $rowVals_376[$i]
}

# Row 377
$dateColStyleSource.Copy()
$ws.Range("A377").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_377 = @(44451,2,1,1,11,3,0,3,3,5,0,2,1,2,0,2,0,2,1,0,5,42,0,0,3,5,0,0,1,0,4,1,4,0,2,0,5,0,10,0,1,127,0,0,0,1,1,0,0,3)
for ($i = 0; $i -lt $rowVals_377.Length; $i++) {
    $ws.Cells.Item(377, $i + 1).Value = $rowVals_377[$i]
}

# Row 378
$dateColStyleSource.Copy()
$ws.Range("A378").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_378 = @(44452,1,0,0,9,7,0,4,1,0,0,0,0,3,0,1,1,2,3,0,4,22,1,1,4,1,0,0,0,0,1,1,0,2,1,0,0,1,4,0,3,83,1,1,0,0,3,0,0,0)
for ($i = 0; $i -lt $rowVals_378.Length; $i++) {
    $ws.Cells.Item(378, $i + 1).Value = $rowVals_378[$i]
}

# Row 379
$dateColStyleSource.Copy()
$ws.Range("A379").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_379 = @(44453,0,1,0,2,0,2,2,0,0,0,1,2,1,0,5,0,0,0,0,5,24,0,0,1,4,1,0,0,0,1,1,4,0,0,0,0,3,7,1,0,70,1,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $rowVals_379.Length; $i++) {
    $ws.Cells.Item(379, $i + 1).Value = $rowVals_379[$i]
}

# Row 380
$dateColStyleSource.Copy()
$ws.Range("A380").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_380 = @(44454,2,0,0,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,16,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,22,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $rowVals_380.Length; $i++) {
    $ws.Cells.Item(380, $i + 1).Value = $rowVals_380[$i]
}

# Row 381
$dateColStyleSource.Copy()
$ws.Range("A381").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_381 = @(44455,0,1,1,10,7,2,3,0,1,0,0,2,9,0,0,0,7,0,0,0,5,0,0,2,2,0,0,1,1,2,0,9,4,1,0,0,2,9,0,4,85,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $rowVals_381.Length; $i++) {
    $ws.Cells.Item(381, $i + 1).Value = $rowVals_381[$i]
}

# Row 382
$dateColStyleSource.Copy()
$ws.Range("A382").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_382 = @(44456,1,1,1,2,1,0,1,2,1,0,5,2,8,0,2,2,2,1,0,4,30,0,0,2,8,0,0,1,0,1,2,2,0,0,0,1,6,1,0,4,97,0,3,0,0,0,0,0,0)
for ($i = 0; $i -lt $rowVals_382.Length; $i++) {
    $ws.Cells.Item(382, $i + 1).Value = $rowVals_382[$i]
}

# Row 383
$dateColStyleSource.Copy()
$ws.Range("A383").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_383 = @(44457,0,1,0,14,5,1,6,0,0,0,4,0,4,0,0,0,3,1,0,0,16,2,1,3,2,0,0,0,3,0,0,2,0,1,0,0,5,6,0,2,83,0,0,0,0,0,1,0,0)
for ($i = 0; $i -lt $rowVals_383.Length; $i++) {
    $ws.Cells.Item(383, $i + 1).Value = $rowVals_383[$i]
}

# Row 384
$dateColStyleSource.Copy()
$ws.Range("A384").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_384 = @(44458,1,1,2,3,1,1,4,1,0,0,2,2,2,0,0,0,4,0,0,3,25,2,1,1,4,0,0,0,1,0,0,1,4,1,2,1,0,9,0,0,81,0,0,2,0,0,0,0,0)
for ($i = 0; $i -lt $rowVals_384.Length; $i++) {
    $ws.Cells.Item(384, $i + 1).Value = $rowVals_384[$i]
}

# Row 385
$dateColStyleSource.Copy()
$ws.Range("A385").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rowVals_385 = @(44459,0,1,0,15,2,3,0,0,1,0,1,1,0,1,1,0,2,0,0,4,12,0,0,3,0,0,0,0,0,4,0,2,0,0,0,3,2,2,0,2,62,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $rowVals_385.Length; $i++) {
    $ws.Cells.Item(385, $i + 1).Value = $rowVals_385[$i]
}
